$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 132; this shifts the existing rows 132-198
# down to 133-199 (dates/volumes/prices move with them), exactly like the
# target diff shows.
$ws.Rows.Item(132).Insert()

# Populate the newly inserted row 132 with its own record (same fixed
# columns as every other row in this block, new date/volume/price data).
$ws.Cells.Item(132, 1).Value = 10
$ws.Cells.Item(132, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(132, 3).Value = "La Araucanía"
$ws.Cells.Item(132, 4).Value = 44839
$ws.Cells.Item(132, 5).Value = 9
$ws.Cells.Item(132, 6).Value = 100114007
$ws.Cells.Item(132, 7).Value = "Jengibre"
$ws.Cells.Item(132, 8).Value = "Sin especificar"
$ws.Cells.Item(132, 9).Value = "Primera"
$ws.Cells.Item(132, 10).Value = 20
$ws.Cells.Item(132, 11).Value = 20000
$ws.Cells.Item(132, 12).Value = 20000
$ws.Cells.Item(132, 13).Value = 20000
$ws.Cells.Item(132, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(132, 15).Value = "Perú"
$ws.Cells.Item(132, 16).Value = 1538
$ws.Cells.Item(132, 17).Value = 13
$ws.Cells.Item(132, 18).Value = "Hortaliza"
